# Update odds values on Sheet1 (Jogos_da_Semana_FlashScore_2024-11-14.xlsx)
# per the committed data refresh ("Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.63
$ws.Range("I2").Value = 2.88
$ws.Range("J2").Value = 3.5
$ws.Range("L2").Value = 3.75
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 26
$ws.Range("AA2").Value = 26
$ws.Range("AH2").Value = 7
$ws.Range("AI2").Value = 13
$ws.Range("AJ2").Value = 12
$ws.Range("AK2").Value = 29
$ws.Range("AL2").Value = 29
$ws.Range("AN2").Value = 4.5
$ws.Range("AO2").Value = 17
$ws.Range("AP2").Value = 29
$ws.Range("AQ2").Value = 51
$ws.Range("AR2").Value = 81
$ws.Range("AW2").Value = 4.75
$ws.Range("AX2").Value = 17
$ws.Range("AY2").Value = 34
$ws.Range("BA2").Value = 101

# Row 3
$ws.Range("I3").Value = 4.33
$ws.Range("J3").Value = 2.75
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 8
$ws.Range("AH3").Value = 9.5
$ws.Range("AI3").Value = 21
$ws.Range("AO3").Value = 11
$ws.Range("AU3").Value = 9.5
$ws.Range("AW3").Value = 6
$ws.Range("AX3").Value = 26
$ws.Range("AY3").Value = 41

# Row 4
$ws.Range("G4").Value = 4.1
$ws.Range("I4").Value = 2.05
$ws.Range("AS4").Value = 351

# Row 5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10

# Row 6
$ws.Range("G6").Value = 1.62
$ws.Range("H6").Value = 3.6
$ws.Range("I6").Value = 5.75
$ws.Range("J6").Value = 2.3
$ws.Range("K6").Value = 2.1
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 3.25
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 1.73
$ws.Range("S6").Value = 1.44
$ws.Range("T6").Value = 2.63
$ws.Range("W6").Value = 6
$ws.Range("AC6").Value = 8.5
$ws.Range("AK6").Value = 67
$ws.Range("AN6").Value = 3.4
$ws.Range("AO6").Value = 8.5
$ws.Range("AT6").Value = 2.63
$ws.Range("AX6").Value = 34

# Row 8
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13

# Row 9
$ws.Range("L9").Value = 7.5
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 8.5
$ws.Range("W9").Value = 5.5
$ws.Range("AC9").Value = 8.5
$ws.Range("AH9").Value = 17
$ws.Range("AV9").Value = 81
$ws.Range("AW9").Value = 8.5

# Row 13
$ws.Range("M13").Value = 1.06
$ws.Range("O13").Value = 1.33
